$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.153961896896362
$ws.Range("B1").Value = 2.536565780639648
$ws.Range("C1").Value = 4.14187479019165
$ws.Range("D1").Value = 3.487235546112061
$ws.Range("E1").Value = 1.217021107673645
